$wb = $excel.ActiveWorkbook

# --- Work on the "Assets" sheet: insert a new row for "Returns Queue" ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Select() | Out-Null

# Insert a new row above row 9 (pushes Pausepoint 5's followers down)
# and populate it with the new "Returns Queue" entry.
$wsAssets.Rows.Item(9).Insert() | Out-Null
$wsAssets.Range("A9").Value = "Returns Queue"
$wsAssets.Range("B9").Value = "Returns Queue"

# Grow Table1 so it covers the newly inserted row too (A1:B20 -> A1:B21).
$lo = $wsAssets.ListObjects.Item(1)
$lo.Resize($wsAssets.Range("A1:B21")) | Out-Null

$wsAssets.Range("C12").Select() | Out-Null

# --- Work on the "Local Config" sheet: just change the selection ---
$wsLocal = $wb.Worksheets.Item("Local Config")
$wsLocal.Select() | Out-Null
$wsLocal.Range("A8:XFD8").Select() | Out-Null

# --- Re-activate Assets as the final active sheet & maximize the window ---
$wsAssets.Select() | Out-Null
$excel.WindowState = -4137
